$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date stamp in C1 of the About sheet, formatted as a (built-in) short date.
# Setting the NumberFormat before the Value keeps the engine from emitting a
# stray custom <numFmt> entry and maps straight onto the built-in numFmtId 14.
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
